$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1064.2142
$ws.Range("I2").Value = 656.8333
$ws.Range("J2").Value = 1369.75
$ws.Range("K2").Value = 656.8333
$ws.Range("L2").Value = 1369.75
$ws.Range("M2").Value = -543.8333
$ws.Range("N2").Value = -1595.75
$ws.Range("H20").Value = 3049.5
$ws.Range("I20").Value = 3049.5
$ws.Range("K20").Value = 3049.5
$ws.Range("M20").Value = -2819.5
$ws.Range("H26").Value = 10844.167
$ws.Range("I26").Value = 5013
$ws.Range("J26").Value = 40000
$ws.Range("K26").Value = 5013
$ws.Range("L26").Value = 40000
$ws.Range("M26").Value = -4669
$ws.Range("N26").Value = -40688
$ws.Range("H35").Value = 3049.5
$ws.Range("I35").Value = 3049.5
$ws.Range("K35").Value = 3049.5
$ws.Range("M35").Value = -2670.5
$ws.Range("H40").Value = 5211.875
$ws.Range("I40").Value = 5115.8335
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 5115.8335
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -4940.8335
$ws.Range("N40").Value = -5850
$ws.Range("H132").Value = 2815843.8
$ws.Range("I132").Value = 2815843.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8447531.399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8445001.399999999
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 1256.421
$ws.Range("I135").Value = 591.5333000000001
$ws.Range("K135").Value = 5323.7997
$ws.Range("M135").Value = -2788.7997
$ws.Range("H137").Value = 13724.186
$ws.Range("I137").Value = 18965.334
$ws.Range("J137").Value = 3241.889
$ws.Range("K137").Value = 56896.00199999999
$ws.Range("L137").Value = 9725.667000000001
$ws.Range("M137").Value = -54346.00199999999
$ws.Range("N137").Value = -14825.667
$ws.Range("H138").Value = 40515.58
$ws.Range("I138").Value = 2008.591
$ws.Range("J138").Value = 252304
$ws.Range("K138").Value = 6025.772999999999
$ws.Range("L138").Value = 756912
$ws.Range("M138").Value = -885.7729999999992
$ws.Range("N138").Value = -767192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17620.842
$ws.Range("I32").Value = 19619.947
$ws.Range("J32").Value = 1628
$ws.Range("K32").Value = 19619.947
$ws.Range("L32").Value = 1628
$ws.Range("M32").Value = -19332.947
$ws.Range("N32").Value = -2202
$ws.Range("H61").Value = 5558.207
$ws.Range("I61").Value = 1149.3889
$ws.Range("J61").Value = 12772.637
$ws.Range("K61").Value = 1149.3889
$ws.Range("L61").Value = 12772.637
$ws.Range("M61").Value = -937.3888999999999
$ws.Range("N61").Value = -13196.637
$ws.Range("H74").Value = 191418.72
$ws.Range("I74").Value = 207568.97
$ws.Range("K74").Value = 207568.97
$ws.Range("M74").Value = -206694.97
$ws.Range("H77").Value = 191418.72
$ws.Range("I77").Value = 207568.97
$ws.Range("K77").Value = 1037844.85
$ws.Range("M77").Value = -1033476.85
$ws.Range("H122").Value = 2105.8333
$ws.Range("I122").Value = 2160.9092
$ws.Range("K122").Value = 6482.7276
$ws.Range("M122").Value = -4032.7276
$ws.Range("H132").Value = 1641.7593
$ws.Range("I132").Value = 1073.8667
$ws.Range("K132").Value = 3221.6001
$ws.Range("M132").Value = -691.6001000000001
$ws.Range("H136").Value = 5558.207
$ws.Range("I136").Value = 1149.3889
$ws.Range("J136").Value = 12772.637
$ws.Range("K136").Value = 3448.1667
$ws.Range("L136").Value = 38317.911
$ws.Range("M136").Value = -898.1666999999998
$ws.Range("N136").Value = -43417.911
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2545.8518
$ws.Range("I134").Value = 1733.7
$ws.Range("J134").Value = 4866.2856
$ws.Range("K134").Value = 5201.1
$ws.Range("L134").Value = 14598.8568
$ws.Range("M134").Value = -2666.1
$ws.Range("N134").Value = -19668.8568
$ws.Range("H140").Value = 97087.55
$ws.Range("J140").Value = 97087.55
$ws.Range("L140").Value = 97087.55
$ws.Range("N140").Value = -107447.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 130.71428
$ws.Range("I7").Value = 138.33333
$ws.Range("K7").Value = 138.33333
$ws.Range("M7").Value = -25.33332999999999
$ws.Range("H16").Value = 2896.2
$ws.Range("J16").Value = 3696.6667
$ws.Range("L16").Value = 3696.6667
$ws.Range("N16").Value = -4270.6667
$ws.Range("H31").Value = 3573424.2
$ws.Range("I31").Value = 4348864
$ws.Range("J31").Value = 6400.4
$ws.Range("K31").Value = 4348864
$ws.Range("L31").Value = 6400.4
$ws.Range("M31").Value = -4348569
$ws.Range("N31").Value = -6990.4
$ws.Range("H34").Value = 3573424.2
$ws.Range("I34").Value = 4348864
$ws.Range("J34").Value = 6400.4
$ws.Range("K34").Value = 4348864
$ws.Range("L34").Value = 6400.4
$ws.Range("M34").Value = -4348662
$ws.Range("N34").Value = -6804.4
$ws.Range("H58").Value = 12805.786
$ws.Range("I58").Value = 1317.1852
$ws.Range("J58").Value = 322998
$ws.Range("K58").Value = 1317.1852
$ws.Range("L58").Value = 322998
$ws.Range("M58").Value = -1114.1852
$ws.Range("N58").Value = -323404
$ws.Range("H94").Value = 1246.3478
$ws.Range("I94").Value = 868.2
$ws.Range("J94").Value = 1537.2307
$ws.Range("K94").Value = 868.2
$ws.Range("L94").Value = 1537.2307
$ws.Range("M94").Value = -417.2
$ws.Range("N94").Value = -2439.2307
$ws.Range("H107").Value = 892.1
$ws.Range("I107").Value = 826.46155
$ws.Range("K107").Value = 826.46155
$ws.Range("M107").Value = 1093.53845
$ws.Range("H113").Value = 2896.2
$ws.Range("J113").Value = 3696.6667
$ws.Range("L113").Value = 3696.6667
$ws.Range("N113").Value = -8036.6667
$ws.Range("H136").Value = 12805.786
$ws.Range("I136").Value = 1317.1852
$ws.Range("J136").Value = 322998
$ws.Range("K136").Value = 3951.5556
$ws.Range("L136").Value = 968994
$ws.Range("M136").Value = -1401.5556
$ws.Range("N136").Value = -974094

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4462.4585
$ws.Range("I62").Value = 1866.6666
$ws.Range("J62").Value = 4833.2856
$ws.Range("K62").Value = 5599.9998
$ws.Range("L62").Value = 14499.8568
$ws.Range("M62").Value = -4913.9998
$ws.Range("N62").Value = -15871.8568
$ws.Range("H65").Value = 4462.4585
$ws.Range("I65").Value = 1866.6666
$ws.Range("J65").Value = 4833.2856
$ws.Range("K65").Value = 16799.9994
$ws.Range("L65").Value = 43499.5704
$ws.Range("M65").Value = -13367.9994
$ws.Range("N65").Value = -50363.5704
$ws.Range("H132").Value = 2230.8572
$ws.Range("J132").Value = 1599
$ws.Range("L132").Value = 14391
$ws.Range("N132").Value = -19451
$ws.Range("H134").Value = 941.2857
$ws.Range("I134").Value = 941.2857
$ws.Range("K134").Value = 2823.8571
$ws.Range("M134").Value = 2246.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 18172.777
$ws.Range("I24").Value = 3000
$ws.Range("J24").Value = 22507.857
$ws.Range("K24").Value = 3000
$ws.Range("L24").Value = 22507.857
$ws.Range("M24").Value = -2827
$ws.Range("N24").Value = -22853.857
$ws.Range("H95").Value = 74564.5
$ws.Range("J95").Value = 74564.5
$ws.Range("L95").Value = 74564.5
$ws.Range("N95").Value = -80056.5
$ws.Range("H107").Value = 419.6
$ws.Range("I107").Value = 166.33333
$ws.Range("K107").Value = 166.33333
$ws.Range("M107").Value = 1753.66667
$ws.Range("H109").Value = 44000
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H132").Value = 2668.182
$ws.Range("I132").Value = 2557.1428
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7671.428400000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5141.428400000001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1376.9286
$ws.Range("I22").Value = 1071.1666
$ws.Range("J22").Value = 1606.25
$ws.Range("K22").Value = 1071.1666
$ws.Range("L22").Value = 1606.25
$ws.Range("M22").Value = -776.1666
$ws.Range("N22").Value = -2196.25
$ws.Range("H27").Value = 1376.9286
$ws.Range("I27").Value = 1071.1666
$ws.Range("J27").Value = 1606.25
$ws.Range("K27").Value = 1071.1666
$ws.Range("L27").Value = 1606.25
$ws.Range("M27").Value = -964.1666
$ws.Range("N27").Value = -1820.25
$ws.Range("H41").Value = 130000
$ws.Range("J41").Value = 130000
$ws.Range("L41").Value = 130000
$ws.Range("N41").Value = -130876
$ws.Range("H46").Value = 4766.769
$ws.Range("I46").Value = 901
$ws.Range("K46").Value = 901
$ws.Range("M46").Value = -713
$ws.Range("H122").Value = 4904.722
$ws.Range("I122").Value = 3803.5386
$ws.Range("J122").Value = 7767.8
$ws.Range("K122").Value = 11410.6158
$ws.Range("L122").Value = 23303.4
$ws.Range("M122").Value = -8960.6158
$ws.Range("N122").Value = -28203.4
$ws.Range("H128").Value = 144974.5
$ws.Range("J128").Value = 144974.5
$ws.Range("L128").Value = 144974.5
$ws.Range("N128").Value = -154934.5
$ws.Range("H136").Value = 2940.795
$ws.Range("I136").Value = 2915.0356
$ws.Range("K136").Value = 8745.106800000001
$ws.Range("M136").Value = -6195.106800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 615.1429000000001
$ws.Range("I113").Value = 353.6316
$ws.Range("K113").Value = 1060.8948
$ws.Range("M113").Value = 1109.1052
$ws.Range("H136").Value = 14623.575
$ws.Range("I136").Value = 15156.421
$ws.Range("J136").Value = 4499.5
$ws.Range("K136").Value = 45469.263
$ws.Range("L136").Value = 13498.5
$ws.Range("M136").Value = -42919.263
$ws.Range("N136").Value = -18598.5

